$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: fix grammar "have to accepted or rejected it" -> "accept or reject"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Proclaims will be attended by members. For every proclaim, members will have to accepted or rejected it. In the first case, a law must be attached. If rejected, a reason must be provided.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Proclaims will be attended by members. For every proclaim, members will have to accept or reject it. In the first case, a law must be attached. If rejected, a reason must be provided.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# Change 4: the stray "_GoBack" bookmark around "Photos are not required..."
# is removed and the two runs it used to separate collapse back together.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$d.Content.Find.Execute(
    "Photos are not required to be stored in the database, but links to external systems like Pin-terest.com or Flickr.com, just to mention a couple of examples.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Photos are not required to be stored in the database, but links to external systems like Pin-terest.com or Flickr.com, just to mention a couple of examples.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: "Create, Update, Delete, List and Show his/her proclaims." is
# shortened to "Manage his/her proclaims." Only the first portion ("Create,
# Update, Delete, List and Show") is replaced so the trailing
# " his/her proclaims." run is left alone.
# ---------------------------------------------------------------------------
$rngManage = $d.Content
$rngManage.Find.Execute("Create, Update, Delete, List and Show", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$manageStart = $rngManage.Start
$manageEnd = $rngManage.End
$manageWord = "Manage"
$insPoint = $d.Range($manageStart, $manageStart)
$insPoint.InsertBefore($manageWord)
$oldPhrase = $d.Range($manageStart + $manageWord.Length, $manageEnd + $manageWord.Length)
$oldPhrase.Delete()

# ---------------------------------------------------------------------------
# Change 3: a new bullet is added right after "An actor who is registered as
# Student must be able to:" (same list level/formatting as the bullet that
# follows it), reading "Create, Update, Delete, List and Show his/her
# proclaims." The document's "_GoBack" bookmark now marks the end of that
# new sentence.
# ---------------------------------------------------------------------------
$rngNext = $d.Content
$rngNext.Find.Execute("Register to the system the different proclaims that he/she may have", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$nextStart = $rngNext.Start
$newSentence = "Create, Update, Delete, List and Show his/her proclaims."
$insPoint2 = $d.Range($nextStart, $nextStart)
$insPoint2.InsertBefore($newSentence + [char]13)

$rngNewSentence = $d.Content
$rngNewSentence.Find.Execute($newSentence, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rngNewSentence.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rngNewSentence) | Out-Null
